$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cell G1 ("Template_name" -> "Template_name_strategy")
$ws.Range("G1").Value = "Template_name_strategy"

# Update the selection to reflect the new active cell/selection (G1 only)
[void]$ws.Range("G1").Select()
